$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new lines (line7, line8) were inserted into the table right after line6, pushing
# the former extr1..extr8 rows down by two. Rows 2-7 (line1..line6) are untouched;
# rows 8-17 are rewritten below to their final (post-insert) contents.
# Columns: A (index), B (name), C (from_bus), D (to_bus), E (in_service)
$rows = @(
    @(8,  6, "line7", 14, 11, $true),
    @(9,  7, "line8", 16, 9,  $false),
    @(10, 8, "extr1", 5,  12, $true),
    @(11, 9, "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $false),
    @(13, 11, "extr4", 7,  8,  $true),
    @(14, 12, "extr5", 9,  11, $true),
    @(15, 13, "extr6", 7,  11, $false),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $true)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# Column A carries the same style as the rest of column A (bold/border/centered) for new rows 16-17 too.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
